# Adds a new study (Mustapich, 2021 - ice hockey) to the study_characteristics
# sheet, and registers "Ice hockey" as a new entry in the sport codebook list
# on the codebook sheet (inserted right before "Dance", renumbering the
# integer codes that follow it).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) study_characteristics sheet: append the new study as row 159.
#    (Written before the codebook edit below so new shared strings land in
#    the same order the original authoring session produced them: title,
#    author, then the "Ice hockey" sport label.)
# ---------------------------------------------------------------------------
$study = $wb.Worksheets.Item("study_characteristics")

$study.Range("A159").Value = 162
$study.Range("D159").Value = "Effects of training load and non-training stress on injury risk in collegiate ice hockey players"
$study.Range("B159").Value = "Mustapich"
$study.Range("C159").Value = 2021
$study.Range("G159").Value = "Ice hockey"
$study.Range("Z159").Value = "Other"
$study.Range("AF159").Value = "Yes"

# ---------------------------------------------------------------------------
# 2) codebook sheet: insert "Ice hockey" into the `sport` code list.
#    Before:  row17 Basketball=0, row18 American Football=10, ... row24 Dance=16,
#             row25 Multiple sports=17
#    After:   row17 Basketball=10, row18 American Football=11, ... row23 Field
#             hockey=16, row24 Ice hockey=17 (NEW), row25 Dance=18,
#             row26 Multiple sports=19
# ---------------------------------------------------------------------------
$codebook = $wb.Worksheets.Item("codebook")

# Insert a fresh row above the current "Dance" row (row 24), pushing Dance,
# Multiple sports, and everything below them down by one row.
$codebook.Rows(24).Insert()

# Renumber the existing codes 0,10,11,12,13,14,15 (Basketball..Field hockey)
# up by one notch so the list stays a clean 1..19 sequence.
$codebook.Range("B17").Value = 10
$codebook.Range("B18").Value = 11
$codebook.Range("B19").Value = 12
$codebook.Range("B20").Value = 13
$codebook.Range("B21").Value = 14
$codebook.Range("B22").Value = 15
$codebook.Range("B23").Value = 16

# Fill in the freshly inserted row with the new "Ice hockey" entry.
$codebook.Range("B24").Value = 17
$codebook.Range("C24").Value = "Ice hockey"

# ---------------------------------------------------------------------------
# 3) View-state touch-ups (active cell / scroll position) matching where the
#    author ended up working.
# ---------------------------------------------------------------------------
$study.Application.ActiveWindow.ScrollRow = 141
$study.Range("T159").Select()

$studyStats = $wb.Worksheets.Item("study_statistical_methods")
$studyStats.Range("C300").Select()

$codebook.Range("B26").Select()

$wb.Save()
